# Apply "poisson_naive pronta para a rodada 27" update to the Fluminense
# historico/rodada_26 sheet.
#
#  - Column A (match id) is renumbered for every data row (2-27).
#  - Column E ("Round") switches from text "Matchweek N" to the bare
#    numeric value N for every data row.
#  - Rows 22 and 23 additionally swap their entire match-stat payload
#    (columns B, C, F, H-BD) with each other; the Opponent column (K)
#    is "Fluminense" in both already, so it is left untouched.
#
# Two of the swapped text cells (B22/B23) hold ISO-style date strings
# ("2023-05-10"/"2023-05-28"). Assigning such a string straight to
# .Value makes Excel's COM layer auto-detect it as a real date (like
# typing it into a cell would), which would store it as a serial
# number with a date NumberFormat instead of the plain text the
# workbook actually uses. Pre-formatting the cell as Text keeps the
# assignment literal; resetting Style back to "Normal" afterwards
# drops the now-unneeded explicit number format again so the cell
# ends up styled exactly like its neighbours.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet



# Row 2
$ws.Cells.Item(2, 1).Value = 3  # A2
$ws.Cells.Item(2, 5).Value = 2  # E2

# Row 3
$ws.Cells.Item(3, 1).Value = 6  # A3
$ws.Cells.Item(3, 5).Value = 4  # E3

# Row 4
$ws.Cells.Item(4, 1).Value = 8  # A4
$ws.Cells.Item(4, 5).Value = 6  # E4

# Row 5
$ws.Cells.Item(5, 1).Value = 12  # A5
$ws.Cells.Item(5, 5).Value = 9  # E5

# Row 6
$ws.Cells.Item(6, 1).Value = 15  # A6
$ws.Cells.Item(6, 5).Value = 11  # E6

# Row 7
$ws.Cells.Item(7, 1).Value = 16  # A7
$ws.Cells.Item(7, 5).Value = 12  # E7

# Row 8
$ws.Cells.Item(8, 1).Value = 19  # A8
$ws.Cells.Item(8, 5).Value = 14  # E8

# Row 9
$ws.Cells.Item(9, 1).Value = 20  # A9
$ws.Cells.Item(9, 5).Value = 15  # E9

# Row 10
$ws.Cells.Item(10, 1).Value = 22  # A10
$ws.Cells.Item(10, 5).Value = 17  # E10

# Row 11
$ws.Cells.Item(11, 1).Value = 24  # A11
$ws.Cells.Item(11, 5).Value = 18  # E11

# Row 12
$ws.Cells.Item(12, 1).Value = 27  # A12
$ws.Cells.Item(12, 5).Value = 20  # E12

# Row 13
$ws.Cells.Item(13, 1).Value = 31  # A13
$ws.Cells.Item(13, 5).Value = 22  # E13

# Row 14
$ws.Cells.Item(14, 1).Value = 33  # A14
$ws.Cells.Item(14, 5).Value = 24  # E14

# Row 15
$ws.Cells.Item(15, 1).Value = 37  # A15
$ws.Cells.Item(15, 5).Value = 26  # E15

# Row 16
$ws.Cells.Item(16, 1).Value = 9  # A16
$ws.Cells.Item(16, 5).Value = 7  # E16

# Row 17
$ws.Cells.Item(17, 1).Value = 17  # A17
$ws.Cells.Item(17, 5).Value = 19  # E17

# Row 18
$ws.Cells.Item(18, 1).Value = 8  # A18
$ws.Cells.Item(18, 5).Value = 3  # E18

# Row 19
$ws.Cells.Item(19, 1).Value = 28  # A19
$ws.Cells.Item(19, 5).Value = 21  # E19

# Row 20
$ws.Cells.Item(20, 1).Value = 18  # A20
$ws.Cells.Item(20, 5).Value = 13  # E20

# Row 21
$ws.Cells.Item(21, 1).Value = 24  # A21
$ws.Cells.Item(21, 5).Value = 25  # E21

# Row 22
$ws.Cells.Item(22, 1).Value = 4  # A22
$ws.Cells.Item(22, 2).NumberFormat = "@"  # B22: force text so the date-like string is not auto-converted
$ws.Cells.Item(22, 2).Value = '2023-05-10'  # B22
$ws.Cells.Item(22, 2).Style = "Normal"  # B22: drop the now-unneeded explicit text format
$ws.Cells.Item(22, 3).Value = '21:30'  # C22
$ws.Cells.Item(22, 5).Value = 5  # E22
$ws.Cells.Item(22, 6).Value = 'Wed'  # F22
$ws.Cells.Item(22, 8).Value = 'L'  # H22
$ws.Cells.Item(22, 9).Value = 0  # I22
$ws.Cells.Item(22, 10).Value = 2  # J22
$ws.Cells.Item(22, 12).Value = 2.1  # L22
$ws.Cells.Item(22, 13).Value = 1.2  # M22
$ws.Cells.Item(22, 14).Value = 48  # N22
$ws.Cells.Item(22, 15).Value = 51846  # O22
$ws.Cells.Item(22, 16).Value = 26  # P22
$ws.Cells.Item(22, 17).Value = 5  # Q22
$ws.Cells.Item(22, 18).Value = 19.2  # R22
$ws.Cells.Item(22, 19).Value = 0  # S22
$ws.Cells.Item(22, 20).Value = 0  # T22
$ws.Cells.Item(22, 23).Value = 1  # W22
$ws.Cells.Item(22, 24).Value = 1.3  # X22
$ws.Cells.Item(22, 25).Value = 0.05  # Y22
$ws.Cells.Item(22, 26).Value = -2.1  # Z22
$ws.Cells.Item(22, 27).Value = -1.3  # AA22
$ws.Cells.Item(22, 28).Value = 6  # AB22
$ws.Cells.Item(22, 29).Value = 4  # AC22
$ws.Cells.Item(22, 30).Value = 66.7  # AD22
$ws.Cells.Item(22, 31).Value = 0  # AE22
$ws.Cells.Item(22, 32).Value = 1.9  # AF22
$ws.Cells.Item(22, 33).Value = -0.1  # AG22
$ws.Cells.Item(22, 34).Value = 7458  # AH22
$ws.Cells.Item(22, 35).Value = 2342  # AI22
$ws.Cells.Item(22, 36).Value = 0  # AJ22
$ws.Cells.Item(22, 37).Value = 1.2  # AK22
$ws.Cells.Item(22, 38).Value = 1.6  # AL22
$ws.Cells.Item(22, 39).Value = 23  # AM22
$ws.Cells.Item(22, 40).Value = 33  # AN22
$ws.Cells.Item(22, 41).Value = 15  # AO22
$ws.Cells.Item(22, 42).Value = 5  # AP22
$ws.Cells.Item(22, 43).Value = 56  # AQ22
$ws.Cells.Item(22, 44).Value = 2  # AR22
$ws.Cells.Item(22, 45).Value = 4  # AS22
$ws.Cells.Item(22, 46).Value = 47  # AT22
$ws.Cells.Item(22, 47).Value = 14  # AU22
$ws.Cells.Item(22, 48).Value = 50  # AV22
$ws.Cells.Item(22, 49).Value = 0  # AW22
$ws.Cells.Item(22, 50).Value = 12  # AX22
$ws.Cells.Item(22, 51).Value = 9  # AY22
$ws.Cells.Item(22, 52).Value = 9  # AZ22
$ws.Cells.Item(22, 53).Value = 5  # BA22
$ws.Cells.Item(22, 54).Value = 9  # BB22
$ws.Cells.Item(22, 56).Value = 'Cruzeiro'  # BD22

# Row 23
$ws.Cells.Item(23, 1).Value = 11  # A23
$ws.Cells.Item(23, 2).NumberFormat = "@"  # B23: force text so the date-like string is not auto-converted
$ws.Cells.Item(23, 2).Value = '2023-05-28'  # B23
$ws.Cells.Item(23, 2).Style = "Normal"  # B23: drop the now-unneeded explicit text format
$ws.Cells.Item(23, 3).Value = '16:00'  # C23
$ws.Cells.Item(23, 5).Value = 8  # E23
$ws.Cells.Item(23, 6).Value = 'Sun'  # F23
$ws.Cells.Item(23, 8).Value = 'W'  # H23
$ws.Cells.Item(23, 9).Value = 2  # I23
$ws.Cells.Item(23, 10).Value = 0  # J23
$ws.Cells.Item(23, 12).Value = 1.4  # L23
$ws.Cells.Item(23, 13).Value = 1.3  # M23
$ws.Cells.Item(23, 14).Value = 30  # N23
$ws.Cells.Item(23, 15).Value = 34624  # O23
$ws.Cells.Item(23, 16).Value = 14  # P23
$ws.Cells.Item(23, 17).Value = 4  # Q23
$ws.Cells.Item(23, 18).Value = 28.6  # R23
$ws.Cells.Item(23, 19).Value = 0.14  # S23
$ws.Cells.Item(23, 20).Value = 0.5  # T23
$ws.Cells.Item(23, 23).Value = 0  # W23
$ws.Cells.Item(23, 24).Value = 1.4  # X23
$ws.Cells.Item(23, 25).Value = 0.1  # Y23
$ws.Cells.Item(23, 26).Value = 0.6  # Z23
$ws.Cells.Item(23, 27).Value = 0.6  # AA23
$ws.Cells.Item(23, 28).Value = 7  # AB23
$ws.Cells.Item(23, 29).Value = 7  # AC23
$ws.Cells.Item(23, 30).Value = 100  # AD23
$ws.Cells.Item(23, 31).Value = 1  # AE23
$ws.Cells.Item(23, 32).Value = 1  # AF23
$ws.Cells.Item(23, 33).Value = 1  # AG23
$ws.Cells.Item(23, 34).Value = 4488  # AH23
$ws.Cells.Item(23, 35).Value = 1700  # AI23
$ws.Cells.Item(23, 36).Value = 2  # AJ23
$ws.Cells.Item(23, 37).Value = 1.3  # AK23
$ws.Cells.Item(23, 38).Value = 0.9  # AL23
$ws.Cells.Item(23, 39).Value = 9  # AM23
$ws.Cells.Item(23, 40).Value = 13  # AN23
$ws.Cells.Item(23, 41).Value = 4  # AO23
$ws.Cells.Item(23, 42).Value = 1  # AP23
$ws.Cells.Item(23, 43).Value = 17  # AQ23
$ws.Cells.Item(23, 44).Value = 1  # AR23
$ws.Cells.Item(23, 45).Value = 3  # AS23
$ws.Cells.Item(23, 46).Value = 17  # AT23
$ws.Cells.Item(23, 47).Value = 6  # AU23
$ws.Cells.Item(23, 48).Value = 20  # AV23
$ws.Cells.Item(23, 49).Value = 4  # AW23
$ws.Cells.Item(23, 50).Value = 11  # AX23
$ws.Cells.Item(23, 51).Value = 12  # AY23
$ws.Cells.Item(23, 52).Value = 2  # AZ23
$ws.Cells.Item(23, 53).Value = 1  # BA23
$ws.Cells.Item(23, 54).Value = 8  # BB23
$ws.Cells.Item(23, 56).Value = 'Corinthians'  # BD23

# Row 24
$ws.Cells.Item(24, 1).Value = 21  # A24
$ws.Cells.Item(24, 5).Value = 23  # E24

# Row 25
$ws.Cells.Item(25, 1).Value = 14  # A25
$ws.Cells.Item(25, 5).Value = 10  # E25

# Row 26
$ws.Cells.Item(26, 1).Value = 15  # A26
$ws.Cells.Item(26, 5).Value = 16  # E26

# Row 27
$ws.Cells.Item(27, 1).Value = 1  # A27
$ws.Cells.Item(27, 5).Value = 1  # E27


Write-Output "Applied rodada_27 poisson_naive update: 148 cells across 26 rows"
